# Update column G ("K" - strikeouts) values for rows 2-26 in the
# houser_adrian save-data sheet. These values were regenerated upstream
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals") and only the G column numbers differ between
# the old and new workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 4
    9  = 0
    10 = 3
    11 = 5
    12 = 3
    13 = 5
    14 = 5
    15 = 1
    16 = 4
    17 = 6
    18 = 3
    19 = 6
    20 = 6
    21 = 3
    22 = 4
    23 = 3
    24 = 4
    25 = 2
    26 = 5
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
